$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F34 was stored as text ("1234567891"); convert it to a real number.
$ws.Range("F34").Value = 1234567891

# New user rows appended below the existing data (rows 35-38).
$ws.Range("A35").Value = "Akash Kawade"
$ws.Range("B35").Value = "Akash5880"
$ws.Range("C35").Value = "Akash@0992"
$ws.Range("D35").Value = "CEO"
$ws.Range("E35").Value = "Male"
$ws.Range("F35").Value = 4567891231
$ws.Range("G35").Value = "Ozar"
$ws.Range("H35").Value = "akash23@gmail.com"

$ws.Range("A36").Value = "Tanmay"
$ws.Range("B36").Value = "Tanmay45"
$ws.Range("C36").Value = "Tanmay@987"
$ws.Range("D36").Value = "HR"
$ws.Range("E36").Value = "Male"
$ws.Range("F36").Value = 1234567897
$ws.Range("G36").Value = "Pune"
$ws.Range("H36").Value = "tanmay@gmail.com"

$ws.Range("A37").Value = "Samarth"
$ws.Range("B37").Value = "Samarth18"
$ws.Range("C37").Value = "Samarth@987"
$ws.Range("D37").Value = "Employee"
$ws.Range("E37").Value = "Male"
$ws.Range("F37").Value = 4557896321
$ws.Range("G37").Value = "Pune"
$ws.Range("H37").Value = "samarth1@gmail.com"

$ws.Range("A38").Value = "Rushi"
$ws.Range("B38").Value = "Rushi45"
$ws.Range("C38").Value = "Rushi@987"
$ws.Range("D38").Value = "HR"
$ws.Range("E38").Value = "Male"
# F38 stays textual ("1234567898"), unlike the other phone numbers above.
# Prefix with an apostrophe to force text entry, then reset the style so
# the cell keeps the plain default formatting (no quote-prefix style).
$ws.Range("F38").Value = "'1234567898"
$ws.Range("F38").Style = "Normal"
$ws.Range("G38").Value = "pune"
$ws.Range("H38").Value = "rushi@gmail.com"
